# Room Database - Processed.xlsx edit script
# Replaces the "Mi 9T" room-database measurements with the "Nexus 5X" ones
# across the Low / Medium / High frequency blocks (and their duplicated
# "All" copies in rows 98:187), then nudges the sheet selection / view to
# match what the author left the workbook in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$lowVals = New-Object 'double[,]' 30,1
$lowVals[0,0] = 111.14243999999999
$lowVals[1,0] = 127.04219999999999
$lowVals[2,0] = 128.13336000000001
$lowVals[3,0] = 127.35396
$lowVals[4,0] = 127.50984
$lowVals[5,0] = 127.50984
$lowVals[6,0] = 126.574559999999
$lowVals[7,0] = 127.50984
$lowVals[8,0] = 126.88632
$lowVals[9,0] = 122.36579999999999
$lowVals[10,0] = 127.8216
$lowVals[11,0] = 127.19808
$lowVals[12,0] = 126.73044
$lowVals[13,0] = 127.04219999999999
$lowVals[14,0] = 127.35396
$lowVals[15,0] = 126.418679999999
$lowVals[16,0] = 127.19808
$lowVals[17,0] = 127.04219999999999
$lowVals[18,0] = 127.04219999999999
$lowVals[19,0] = 127.04219999999999
$lowVals[20,0] = 127.19808
$lowVals[21,0] = 127.97748
$lowVals[22,0] = 127.04219999999999
$lowVals[23,0] = 126.88632
$lowVals[24,0] = 127.35396
$lowVals[25,0] = 127.19808
$lowVals[26,0] = 121.898159999999
$lowVals[27,0] = 126.73044
$lowVals[28,0] = 127.35396
$lowVals[29,0] = 127.19808

$medVals = New-Object 'double[,]' 30,1
$medVals[0,0] = 113.247828
$medVals[1,0] = 132.04414800000001
$medVals[2,0] = 115.91064
$medVals[3,0] = 116.067275999999
$medVals[4,0] = 133.92377999999999
$medVals[5,0] = 111.838104
$medVals[6,0] = 116.067275999999
$medVals[7,0] = 116.067275999999
$medVals[8,0] = 116.537184
$medVals[9,0] = 133.45387199999999
$medVals[10,0] = 116.380547999999
$medVals[11,0] = 115.440732
$medVals[12,0] = 133.297236
$medVals[13,0] = 133.297236
$medVals[14,0] = 133.76714399999901
$medVals[15,0] = 116.067275999999
$medVals[16,0] = 114.97082399999999
$medVals[17,0] = 115.28409600000001
$medVals[18,0] = 115.597368
$medVals[19,0] = 118.260179999999
$medVals[20,0] = 130.94769600000001
$medVals[21,0] = 131.73087599999999
$medVals[22,0] = 124.99552799999999
$medVals[23,0] = 132.35741999999999
$medVals[24,0] = 131.41760399999899
$medVals[25,0] = 131.26096799999999
$medVals[26,0] = 131.57424
$medVals[27,0] = 131.41760399999899
$medVals[28,0] = 131.88751199999999
$medVals[29,0] = 134.863596

$highVals = New-Object 'double[,]' 30,1
$highVals[0,0] = 113.5611
$highVals[1,0] = 129.537972
$highVals[2,0] = 111.211559999999
$highVals[3,0] = 128.911428
$highVals[4,0] = 111.211559999999
$highVals[5,0] = 128.911428
$highVals[6,0] = 128.911428
$highVals[7,0] = 129.22469999999899
$highVals[8,0] = 129.06806399999999
$highVals[9,0] = 128.59815599999999
$highVals[10,0] = 129.85124399999901
$highVals[11,0] = 128.28488399999901
$highVals[12,0] = 129.06806399999999
$highVals[13,0] = 128.75479200000001
$highVals[14,0] = 129.06806399999999
$highVals[15,0] = 128.75479200000001
$highVals[16,0] = 128.75479200000001
$highVals[17,0] = 129.06806399999999
$highVals[18,0] = 129.06806399999999
$highVals[19,0] = 128.75479200000001
$highVals[20,0] = 128.75479200000001
$highVals[21,0] = 128.59815599999999
$highVals[22,0] = 124.99552799999999
$highVals[23,0] = 129.22469999999899
$highVals[24,0] = 128.911428
$highVals[25,0] = 128.59815599999999
$highVals[26,0] = 129.381336
$highVals[27,0] = 128.75479200000001
$highVals[28,0] = 128.911428
$highVals[29,0] = 128.59815599999999

$allVals = New-Object 'double[,]' 90,1
$allVals[0,0] = 111.14243999999999
$allVals[1,0] = 127.04219999999999
$allVals[2,0] = 128.13336000000001
$allVals[3,0] = 127.35396
$allVals[4,0] = 127.50984
$allVals[5,0] = 127.50984
$allVals[6,0] = 126.574559999999
$allVals[7,0] = 127.50984
$allVals[8,0] = 126.88632
$allVals[9,0] = 122.36579999999999
$allVals[10,0] = 127.8216
$allVals[11,0] = 127.19808
$allVals[12,0] = 126.73044
$allVals[13,0] = 127.04219999999999
$allVals[14,0] = 127.35396
$allVals[15,0] = 126.418679999999
$allVals[16,0] = 127.19808
$allVals[17,0] = 127.04219999999999
$allVals[18,0] = 127.04219999999999
$allVals[19,0] = 127.04219999999999
$allVals[20,0] = 127.19808
$allVals[21,0] = 127.97748
$allVals[22,0] = 127.04219999999999
$allVals[23,0] = 126.88632
$allVals[24,0] = 127.35396
$allVals[25,0] = 127.19808
$allVals[26,0] = 121.898159999999
$allVals[27,0] = 126.73044
$allVals[28,0] = 127.35396
$allVals[29,0] = 127.19808
$allVals[30,0] = 113.247828
$allVals[31,0] = 132.04414800000001
$allVals[32,0] = 115.91064
$allVals[33,0] = 116.067275999999
$allVals[34,0] = 133.92377999999999
$allVals[35,0] = 111.838104
$allVals[36,0] = 116.067275999999
$allVals[37,0] = 116.067275999999
$allVals[38,0] = 116.537184
$allVals[39,0] = 133.45387199999999
$allVals[40,0] = 116.380547999999
$allVals[41,0] = 115.440732
$allVals[42,0] = 133.297236
$allVals[43,0] = 133.297236
$allVals[44,0] = 133.76714399999901
$allVals[45,0] = 116.067275999999
$allVals[46,0] = 114.97082399999999
$allVals[47,0] = 115.28409600000001
$allVals[48,0] = 115.597368
$allVals[49,0] = 118.260179999999
$allVals[50,0] = 130.94769600000001
$allVals[51,0] = 131.73087599999999
$allVals[52,0] = 124.99552799999999
$allVals[53,0] = 132.35741999999999
$allVals[54,0] = 131.41760399999899
$allVals[55,0] = 131.26096799999999
$allVals[56,0] = 131.57424
$allVals[57,0] = 131.41760399999899
$allVals[58,0] = 131.88751199999999
$allVals[59,0] = 134.863596
$allVals[60,0] = 113.5611
$allVals[61,0] = 129.537972
$allVals[62,0] = 111.211559999999
$allVals[63,0] = 128.911428
$allVals[64,0] = 111.211559999999
$allVals[65,0] = 128.911428
$allVals[66,0] = 128.911428
$allVals[67,0] = 129.22469999999899
$allVals[68,0] = 129.06806399999999
$allVals[69,0] = 128.59815599999999
$allVals[70,0] = 129.85124399999901
$allVals[71,0] = 128.28488399999901
$allVals[72,0] = 129.06806399999999
$allVals[73,0] = 128.75479200000001
$allVals[74,0] = 129.06806399999999
$allVals[75,0] = 128.75479200000001
$allVals[76,0] = 128.75479200000001
$allVals[77,0] = 129.06806399999999
$allVals[78,0] = 129.06806399999999
$allVals[79,0] = 128.75479200000001
$allVals[80,0] = 128.75479200000001
$allVals[81,0] = 128.59815599999999
$allVals[82,0] = 124.99552799999999
$allVals[83,0] = 129.22469999999899
$allVals[84,0] = 128.911428
$allVals[85,0] = 128.59815599999999
$allVals[86,0] = 129.381336
$allVals[87,0] = 128.75479200000001
$allVals[88,0] = 128.911428
$allVals[89,0] = 128.59815599999999

# --- Write the new measurement values -------------------------------------
$ws.Range("B2:B31").Value2   = $lowVals
$ws.Range("B34:B63").Value2  = $medVals
$ws.Range("B66:B95").Value2  = $highVals
$ws.Range("B98:B187").Value2 = $allVals

# --- Force a full recalculation so the AVERAGE/MIN/MAX/QUARTILE/STDEV -----
# --- helper formulas in columns D:E pick up the new numbers ---------------
$excel.CalculateFullRebuild()

# --- Best-effort cosmetic state: scroll/selection left by the author ------
$ws.Range("B99:B187").Select()
